$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.827.35'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '2.546.76'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.45'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.48'
$ws.Range('E6').Value = '  -1.33%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.523'
$ws.Range('E8').Value = '  +2.05%  '
$ws.Range('D9').Value = '2.544.75'
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.164'
$ws.Range('E10').Value = '  -2.20%  '
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.359'
$ws.Range('E12').Value = '  +2.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.93'
$ws.Range('E13').Value = '  +1.43%  '
$ws.Range('D14').Value = '3.014.52'
$ws.Range('E14').Value = '  -0.93%  '
$ws.Range('D15').Value = '70.708.07'
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000179'
$ws.Range('E16').Value = '  -3.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.34'
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').Value = '2.566.92'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.97'
$ws.Range('E19').Value = '  +3.56%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.53'
$ws.Range('E20').Value = '  -2.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '353.62'
$ws.Range('E21').Value = '  -3.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.95'
$ws.Range('E22').Value = '  -1.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.03'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.44'
$ws.Range('E25').Value = '  +0.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.05'
$ws.Range('E26').Value = '  -2.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.98'
$ws.Range('E27').Value = '  -3.98%  '
$ws.Range('D28').Value = '2.706.92'
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').Value = '0.0₃0916'
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.95'
$ws.Range('E31').Value = '  +1.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '468.79'
$ws.Range('E32').Value = '  -3.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.27'
$ws.Range('E33').Value = '  -3.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.76'
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('E36').Value = '  +2.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '156.81'
$ws.Range('E37').Value = '  -0.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.86'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.10'
$ws.Range('E39').Value = '  +1.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.85'
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.321'
$ws.Range('E42').Value = '  -0.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.62'
$ws.Range('E43').Value = '  -4.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.34'
$ws.Range('E44').Value = '  -6.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.16'
$ws.Range('E45').Value = '  -13.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.62'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '144.38'
$ws.Range('E47').Value = '  -1.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.534'
$ws.Range('E48').Value = '  -0.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.53'
$ws.Range('E49').Value = '  -1.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.61'
$ws.Range('E50').Value = '  -2.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0736'
$ws.Range('E51').Value = '  -0.01%  '
